$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("R2").Value = 1.77390151467452
$ws.Range("R3").Value = 1.733144255704353
$ws.Range("R4").Value = 39.27
$ws.Range("R5").Value = 3.953271223769448
$ws.Range("R6").Value = 34.44
$ws.Range("R7").Value = 41.46613151759277
$ws.Range("R8").Value = 41.43457877430134
$ws.Range("R9").Value = 71.92
$ws.Range("R11").Value = 11.61185150262494
$ws.Range("R12").Value = 21.22
